$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue 2 4 "29.234.08"
$ws.Cells.Item(2, 5).Value = "  +0.44%  "

# Row 3 - Ethereum
Set-TextValue 3 4 "1.857.29"
$ws.Cells.Item(3, 5).Value = "  +0.40%  "

# Row 4 - TetherUSD
$ws.Cells.Item(4, 5).Value = "  -0.01%  "

# Row 5 - XRP
Set-TextValue 5 4 "0.7085"
$ws.Cells.Item(5, 5).Value = "  +2.07%  "

# Row 6 - BNB
Set-TextValue 6 4 "238.44"
$ws.Cells.Item(6, 5).Value = "  +0.21%  "

# Row 7 - USDC
$ws.Cells.Item(7, 5).Value = "  -0.01%  "

# Row 8 - Dogecoin
Set-TextValue 8 4 "0.07995"
$ws.Cells.Item(8, 5).Value = "  +3.80%  "

# Row 9 - Cardano
Set-TextValue 9 4 "0.3030"
$ws.Cells.Item(9, 5).Value = "  -0.16%  "

# Row 10 - Solana
Set-TextValue 10 4 "23.48"
$ws.Cells.Item(10, 5).Value = "  +1.00%  "

# Row 11 - TRON
Set-TextValue 11 4 "0.08198"
$ws.Cells.Item(11, 5).Value = "  +1.06%  "

# Row 12 - Polkadot
Set-TextValue 12 4 "5.182"
$ws.Cells.Item(12, 5).Value = "  -0.56%  "

# Row 13 - Polygon
Set-TextValue 13 4 "0.7030"
$ws.Cells.Item(13, 5).Value = "  -3.14%  "

# Row 14 - now WrappedEther (was Litecoin)
$ws.Cells.Item(14, 2).Value = "WrappedEther"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue 14 4 "1.773.65"
$ws.Cells.Item(14, 5).Value = "  -3.93%  "

# Row 15 - now Litecoin (was WrappedEther)
$ws.Cells.Item(15, 2).Value = "Litecoin"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue 15 4 "89.62"
$ws.Cells.Item(15, 5).Value = "  +0.67%  "

# Row 16 - WrappedBTC
Set-TextValue 16 4 "29.113.11"
$ws.Cells.Item(16, 5).Value = "  +0.06%  "

# Row 17 - Uniswap
Set-TextValue 17 4 "5.806"
$ws.Cells.Item(17, 5).Value = "  +0.98%  "

# Row 18 - ShibaInu
Set-TextValue 18 4 "0.000007886"
$ws.Cells.Item(18, 5).Value = "  +1.48%  "

# Row 19 - Avalanche
Set-TextValue 19 4 "13.23"
$ws.Cells.Item(19, 5).Value = "  +0.18%  "

# Row 20 - BitcoinCash
Set-TextValue 20 4 "237.61"
$ws.Cells.Item(20, 5).Value = "  +0.60%  "

# Row 21 - Dai
Set-TextValue 21 4 "0.9986"
$ws.Cells.Item(21, 5).Value = "  -0.13%  "

# Row 22 - BinanceUSD
Set-TextValue 22 4 "1.000"
$ws.Cells.Item(22, 5).Value = "  -0.06%  "

# Row 23 - WrappedliquidstakedEther2.0
Set-TextValue 23 4 "2.037.46"
$ws.Cells.Item(23, 5).Value = "  -2.02%  "

# Row 24 - Chainlink
Set-TextValue 24 4 "7.461"
$ws.Cells.Item(24, 5).Value = "  -1.76%  "

# Row 25 - Monero
Set-TextValue 25 4 "162.87"
$ws.Cells.Item(25, 5).Value = "  +1.06%  "

# Row 26 - Cosmos
Set-TextValue 26 4 "8.895"
$ws.Cells.Item(26, 5).Value = "  -0.94%  "

# Row 27 - Stellar
Set-TextValue 27 4 "0.1437"
$ws.Cells.Item(27, 5).Value = "  +0.29%  "

# Row 28 - EthereumClassic
$ws.Cells.Item(28, 5).Value = "  +0.38%  "

# Row 29 - LidoDAOToken
Set-TextValue 29 4 "1.919"
$ws.Cells.Item(29, 5).Value = "  -3.24%  "

# Row 30 - Toncoin
Set-TextValue 30 4 "1.421"
$ws.Cells.Item(30, 5).Value = "  +1.57%  "

# Row 31 - PancakeSwap
Set-TextValue 31 4 "1.478"
$ws.Cells.Item(31, 5).Value = "  -0.42%  "

# Row 32 - Filecoin
Set-TextValue 32 4 "4.357"
$ws.Cells.Item(32, 5).Value = "  -3.05%  "

# Row 33 - InternetComputer(DFINITY)
Set-TextValue 33 4 "4.021"
$ws.Cells.Item(33, 5).Value = "  +0.12%  "

# Row 34 - Hedera
Set-TextValue 34 4 "0.05192"
$ws.Cells.Item(34, 5).Value = "  -0.67%  "

# Row 35 - ARBITRUM (only price changes, E unchanged)
Set-TextValue 35 4 "1.159"

# Row 36 - ImmutableX
Set-TextValue 36 4 "0.7131"
$ws.Cells.Item(36, 5).Value = "  +1.85%  "

# Row 37 - Frax
Set-TextValue 37 4 "0.9973"
$ws.Cells.Item(37, 5).Value = "  -2.51%  "

# Row 38 - HuobiToken
Set-TextValue 38 4 "2.668"
$ws.Cells.Item(38, 5).Value = "  +0.41%  "

# Row 39 - VeChain
Set-TextValue 39 4 "0.01851"
$ws.Cells.Item(39, 5).Value = "  -0.04%  "

# Row 40 - MXToken
$ws.Cells.Item(40, 5).Value = "  +1.57%  "

# Row 41 - TrustWalletToken
Set-TextValue 41 4 "0.9327"
$ws.Cells.Item(41, 5).Value = "  +1.45%  "

# Row 42 - Maker
Set-TextValue 42 4 "1.138.75"
$ws.Cells.Item(42, 5).Value = "  +5.26%  "

# Row 43 - FraxShare
Set-TextValue 43 4 "5.917"
$ws.Cells.Item(43, 5).Value = "  -1.57%  "

# Row 44 - TheSandbox
Set-TextValue 44 4 "0.4257"
$ws.Cells.Item(44, 5).Value = "  -0.02%  "

# Row 45 - Aave
Set-TextValue 45 4 "70.20"
$ws.Cells.Item(45, 5).Value = "  -0.15%  "

# Row 46 - PaxDollar
$ws.Cells.Item(46, 5).Value = "  -0.06%  "

# Row 47 - Quant
Set-TextValue 47 4 "102.42"
$ws.Cells.Item(47, 5).Value = "  -0.70%  "

# Row 48 - Mantle
Set-TextValue 48 4 "0.5332"
$ws.Cells.Item(48, 5).Value = "  -4.40%  "

# Row 49 - RenderToken
$ws.Cells.Item(49, 5).Value = "  -0.62%  "

# Row 50 - EnergySwap
$ws.Cells.Item(50, 5).Value = "  +0.47%  "

# Row 51 - Aptos
Set-TextValue 51 4 "6.949"
$ws.Cells.Item(51, 5).Value = "  -0.79%  "
